$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8334118723869324
$ws.Range("B1").Value = 2.043908596038818
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.611337423324585
$ws.Range("E1").Value = 0.4790465831756592
